$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data was shifted one column to the left (old column A, which held the
# per-row "11-texon"/"14-texon" gene counts, was removed). Deleting column A
# with a left-shift reproduces this exactly: every remaining cell (values,
# shared-string refs, and styles) slides into the previous column, including
# the bold/border/center header style moving onto the new A1.
$ws.Columns("A:A").Delete(-4159)
